$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 45019
$ws.Range("C3").Value = 45020
$ws.Range("C4").Value = 45021
$ws.Range("C5").Value = 45022
$ws.Range("C6").Value = 45023
$ws.Range("C7").Value = 45024
